$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Wrap the existing QUERY formula with INDEX(...,1,1) for the whole A1:A100 range.
# A1 holds a standalone formula; A2:A65 and A66:A100 are two shared-formula groups.
$ws.Range("A1").Formula = '=INDEX(QUERY("ShuffledDataSet100","ShuffledDataSet100_100_local"),1,1)'
$ws.Range("A2:A65").Formula = '=INDEX(QUERY("ShuffledDataSet100","ShuffledDataSet100_100_local"),1,1)'
$ws.Range("A66:A100").Formula = '=INDEX(QUERY("ShuffledDataSet100","ShuffledDataSet100_100_local"),1,1)'

# Widen column A so its stored width tracks the longer formula text
# (bestFit recompute for "=INDEX(QUERY(...),1,1)" -> ~34.71 characters).
$ws.Columns.Item(1).ColumnWidth = 33.8

# Collapse the old whole-row selection (A1:XFD1) down to a single cell.
$ws.Range("A1").Select()
